# Updates the cryptos list with the latest price/volume data.
# - Columns D (Price) and E (Volume(1h)) get refreshed text values for
#   most rows.
# - Rows 48 and 49 (EnergySwap / Mantle) swap places along with their
#   updated price/volume figures.
#
# Column D values are forced to remain plain text (matching the original
# inlineStr cells) even when they look like numbers, by temporarily
# applying a text number format before assignment and then restoring the
# cell style so no extra formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.807.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.637.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0641"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.89"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0786"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.648.00"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("E13").Value = "  -0.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.861.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("E16").Value = "  +2.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.816.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("E20").Value = "  +2.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "194.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("E23").Value = "  +0.93%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("E25").Value = "  -2.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("E27").Value = "  -4.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.83%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  +1.44%  "
$ws.Range("E32").Value = "  +1.11%  "
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("E34").Value = "  +2.06%  "
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("E36").Value = "  -0.59%  "
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.112.36"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.72%  "
$ws.Range("E40").Value = "  +0.38%  "
$ws.Range("E41").Value = "  +0.43%  "
$ws.Range("E42").Value = "  +1.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.24"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.799"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("E45").Value = "  -1.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +11.86%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.418"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.79%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.23%  "
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.51%  "
